$wb = $excel.ActiveWorkbook

# Data to insert into the new "fisherman_days" column (column C) for each
# sheet, keyed by sheet/year. Values read off row 2..N in the order the
# rows already appear in each sheet.
$fishermanDays = @{
    "1938" = @(6250, 7323, 11119, 11407, 12658, 26802, 22214, 18822, 6106, 950)
    "1939" = @(11398, 17300, 28016, 41022, 17722, 6597, 2908)
    "1940" = @(14314, 12481, 17497, 28944, 35597, 29960, 11522, 3171)
    "1941" = @(15438, 13319, 17878, 20094, 34216, 20101, 5824, 1538)
    "1942" = @(109, 1157, 502, 7189, 16200, 10088, 5281, 522)
}

foreach ($sheetName in @("1938", "1939", "1940", "1941", "1942")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Insert a new blank column in front of the existing "rock_cod" data
    # column (column C), shifting rock_cod / NA-or-count / total_boathouses
    # one column to the right.
    $ws.Columns.Item(3).Insert()

    # Keep the new column's width consistent with its neighbours (13.5),
    # matching the sheet's existing formatting convention.
    $ws.Columns.Item(3).ColumnWidth = 12.71

    # Header + values for the newly-opened column.
    $ws.Range("C1").Value = "fisherman_days"

    $values = $fishermanDays[$sheetName]
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($i + 2, 3).Value = $values[$i]
    }
}

# Move the active tab / selection: 1940 (sheet index 3) becomes the
# selected tab, replacing 1942 (sheet index 5).
$wb.Worksheets.Item("1938").Range("E23").Select()
$wb.Worksheets.Item("1940").Range("C18").Select()
$wb.Worksheets.Item("1941").Range("C10").Select()
$wb.Worksheets.Item("1942").Range("F21").Select()

$wb.Worksheets.Item("1940").Activate()
